$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("f1_score")
$ws.Range("B4").Value = "0.518 (0.462 Â± 0.033)"
$ws.Range("C4").Value = "0.639 (0.567 Â± 0.037)"
$ws.Range("D4").Value = "0.711 (0.677 Â± 0.016)"
$ws.Range("E4").Value = "0.758 (0.728 Â± 0.020)"
$ws.Range("F4").Value = "0.248 (0.234 Â± 0.009)"

$ws.Range("B6").Value = "0.739 (0.690 Â± 0.026)"
$ws.Range("C6").Value = "0.658 (0.625 Â± 0.020)"
$ws.Range("D6").Value = "0.807 (0.775 Â± 0.015)"
$ws.Range("E6").Value = "0.816 (0.791 Â± 0.016)"
$ws.Range("F6").Value = "0.630 (0.597 Â± 0.019)"

$ws.Range("B14").Value = "N/A"
$ws.Range("C14").Value = "N/A"
$ws.Range("D14").Value = "N/A"
$ws.Range("E14").Value = "N/A"
$ws.Range("F14").Value = "N/A"


$ws = $wb.Worksheets.Item("training_time")
$ws.Range("B4").Value = "00:00:27 (00:00:38 Â± 00:00:10)"
$ws.Range("C4").Value = "00:01:48 (00:02:20 Â± 00:00:32)"
$ws.Range("D4").Value = "00:03:15 (00:03:49 Â± 00:00:36)"
$ws.Range("E4").Value = "00:01:07 (00:01:27 Â± 00:00:19)"
$ws.Range("F4").Value = "00:00:52 (00:01:07 Â± 00:00:15)"

$ws.Range("B6").Value = "00:04:56 (00:05:00 Â± 00:00:02)"
$ws.Range("C6").Value = "00:04:57 (00:05:01 Â± 00:00:01)"
$ws.Range("D6").Value = "00:04:56 (00:05:01 Â± 00:00:02)"
$ws.Range("E6").Value = "00:04:56 (00:05:00 Â± 00:00:02)"
$ws.Range("F6").Value = "00:04:54 (00:05:01 Â± 00:00:03)"

$ws.Range("B14").Value = "N/A"
$ws.Range("C14").Value = "N/A"
$ws.Range("D14").Value = "N/A"
$ws.Range("E14").Value = "N/A"
$ws.Range("F14").Value = "N/A"


$ws = $wb.Worksheets.Item("test_time")
$ws.Range("B4").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("C4").Value = "00:00:00 (00:00:00 Â± 00:00:00)"
$ws.Range("D4").Value = "00:00:10 (00:00:10 Â± 00:00:00)"
$ws.Range("E4").Value = "00:00:01 (00:00:01 Â± 00:00:00)"
$ws.Range("F4").Value = "00:00:00 (00:00:00 Â± 00:00:00)"

$ws.Range("B6").Value = "00:00:00 (00:00:02 Â± 00:00:01)"
$ws.Range("C6").Value = "00:00:00 (00:00:02 Â± 00:00:01)"
$ws.Range("D6").Value = "00:00:00 (00:00:02 Â± 00:00:01)"
$ws.Range("E6").Value = "00:00:00 (00:00:04 Â± 00:00:02)"
$ws.Range("F6").Value = "00:00:00 (00:00:05 Â± 00:00:06)"

$ws.Range("B14").Value = "N/A"
$ws.Range("C14").Value = "N/A"
$ws.Range("D14").Value = "N/A"
$ws.Range("E14").Value = "N/A"
$ws.Range("F14").Value = "N/A"


$ws = $wb.Worksheets.Item("missing_runs")
$ws.Range("B4").Value = "[]"
$ws.Range("C4").Value = "[]"
$ws.Range("D4").Value = "[]"
$ws.Range("E4").Value = "[]"
$ws.Range("F4").Value = "[]"

$ws.Range("B6").Value = "[]"
$ws.Range("C6").Value = "[]"
$ws.Range("D6").Value = "[]"
$ws.Range("E6").Value = "[]"
$ws.Range("F6").Value = "[]"

$ws.Range("B14").Value = "N/A"
$ws.Range("C14").Value = "N/A"
$ws.Range("D14").Value = "N/A"
$ws.Range("E14").Value = "N/A"
$ws.Range("F14").Value = "N/A"


$ws = $wb.Worksheets.Item("best_seed")
$ws.Range("B4").Value = 19
$ws.Range("C4").Value = 2
$ws.Range("D4").Value = 43
$ws.Range("E4").Value = 43
$ws.Range("F4").Value = 3

$ws.Range("B6").Value = 71
$ws.Range("C6").Value = 67
$ws.Range("D6").Value = 19
$ws.Range("E6").Value = 19
$ws.Range("F6").Value = 13

$ws.Range("B14").Value = "N/A"
$ws.Range("C14").Value = "N/A"
$ws.Range("D14").Value = "N/A"
$ws.Range("E14").Value = "N/A"
$ws.Range("F14").Value = "N/A"

